# Rename the sheet from "example_384_well" to "Sheet1" and update the
# active selection on that sheet from N46:Y53 (anchor N46) to just E32,
# matching a resave of the workbook from Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"
$ws.Range("E32").Select()
